# Update "想去人数" (interest count) figures in the F column across the
# relevant worksheets, as published with the gh-pages output generated
# at 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 336
$ws1.Range("F11").Value = 1030
$ws1.Range("F13").Value = 1734
$ws1.Range("F15").Value = 6172
$ws1.Range("F21").Value = 4833
$ws1.Range("F23").Value = 359
$ws1.Range("F31").Value = 1051
$ws1.Range("F33").Value = 100
$ws1.Range("F34").Value = 100

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F24").Value = 6536

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 1296
$ws3.Range("F10").Value = 8883
$ws3.Range("F11").Value = 973
$ws3.Range("F12").Value = 77

# 全部类型 (All types, merged view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 336
$ws4.Range("F7").Value  = 973
$ws4.Range("F9").Value  = 77
$ws4.Range("F13").Value = 1030
$ws4.Range("F17").Value = 6172
$ws4.Range("F22").Value = 4833
$ws4.Range("F24").Value = 359
$ws4.Range("F28").Value = 1051
$ws4.Range("F30").Value = 100
$ws4.Range("F31").Value = 100
